{"js": "// Remove the three \"papers read\" bullet paragraphs under\n// \"Point Clouds/Collision detection:\" (the ones about \"Collision and\n// proximity queries...\", \"...pproximating polyhedra...\", and \"Point cloud\n// collision detection, in Eurographics (2004)\"), and move the `_GoBack`\n// bookmark from the end of the \"Point Clouds/Collision detection:\" heading\n// paragraph to the start of the next remaining bullet (\"Efficient bounds\n// for point-based animations...\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst textsToRemove = [\n  \"Collision and proximity queries, in Handbook of Discrete and Computational Geometry\",\n  \"pproximating polyhedra with spheres for time-critical collision detection.\",\n  \"Point cloud collision detection, in Eurographics (2004)\",\n];\n\nconst toDelete = [];\nlet nextParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (textsToRemove.indexOf(text) !== -1) {\n    toDelete.push(para);\n  }\n}\n\n// Delete the three paragraphs.\nfor (const para of toDelete) {\n  para.delete();\n}\nawait context.sync();\n\n// Drop the old `_GoBack` bookmark (currently sitting right after\n// \"Point Clouds/Collision detection:\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-find paragraphs now that the three bullets are gone, and locate the\n// bullet that follows them (\"Efficient bounds for point-based animations...\").\nconst remaining = context.document.body.paragraphs;\nremaining.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < remaining.items.length; i++) {\n  if (remaining.items[i].text.trim().indexOf(\"Efficient bounds for point-based animations\") === 0) {\n    nextParagraph = remaining.items[i];\n    break;\n  }\n}\n\nif (nextParagraph) {\n  const startRange = nextParagraph.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Remove the three \"papers read\" bullet paragraphs under\n# \"Point Clouds/Collision detection:\" (the ones about \"Collision and\n# proximity queries...\", \"...pproximating polyhedra...\", and \"Point cloud\n# collision detection, in Eurographics (2004)\"), and move the `_GoBack`\n# bookmark from the end of the \"Point Clouds/Collision detection:\" heading\n# paragraph to the start of the next remaining bullet (\"Efficient bounds\n# for point-based animations...\").\n\n$d = $word.ActiveDocument\n\n$toDeleteTexts = @(\n  \"Collision and proximity queries, in Handbook of Discrete and Computational Geometry\",\n  \"pproximating polyhedra with spheres for time-critical collision detection.\",\n  \"Point cloud collision detection, in Eurographics (2004)\"\n)\n\n# Snapshot the paragraphs first, then delete the matching ones back-to-front\n# so the still-pending ranges don't shift out from under us.\n$paras = @()\nforeach ($p in $d.Paragraphs) {\n    $paras += $p\n}\nfor ($i = $paras.Count - 1; $i -ge 0; $i--) {\n    $p = $paras[$i]\n    if ($toDeleteTexts -contains $p.Range.Text.Trim()) {\n        $p.Range.Delete()\n    }\n}\n\n# Drop the old `_GoBack` bookmark (currently sitting right after\n# \"Point Clouds/Collision detection:\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-add it, collapsed, at the start of the \"Efficient bounds...\" paragraph.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim().StartsWith(\"Efficient bounds for point-based animations\")) {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $r = $target.Range.Duplicate\n    $r.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $r)\n}\n"}
